$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B13").Value = "2022/3/26完成"

$ws.Range("B13").Select()
